# Actualización desde MV -datos-
# Adds a new daily row (07-09-2021) to the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date label as literal text (not an auto-converted date serial).
# Force a text number format first so Excel keeps the "07-09-2021" string,
# then restore the cell's default style so no stray formatting is left
# behind on the new row.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "07-09-2021"
$ws.Range("A19").Style = "Normal"

$ws.Range("B19").Value = 10000
$ws.Range("D19").Value = 0
